$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data
$ws.Range("D2").Value = '68.423.48'
$ws.Range("E2").Value = '  +1.60%  '

$ws.Range("D3").Value = '3.742.76'
$ws.Range("E3").Value = '  -0.28%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.01%  '

$ws.Range("D7").Value = '3.741.45'
$ws.Range("E7").Value = '  -0.31%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.519'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.78%  '

$ws.Range("E10").Value = '  -3.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.47'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.22%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.449'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.58%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000259'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.79%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.32'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.18%  '

$ws.Range("D15").Value = '4.371.76'
$ws.Range("E15").Value = '  -0.24%  '

$ws.Range("D16").Value = '3.742.69'
$ws.Range("E16").Value = '  -0.38%  '

$ws.Range("D17").Value = '68.425.99'
$ws.Range("E17").Value = '  +1.74%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.94'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.34%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.02'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.09%  '

$ws.Range("E20").Value = '  -0.24%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.70'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '467.17'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.700'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.50%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.62%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000145'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.41%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.19'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.38%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.44%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.13'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.77%  '

$ws.Range("E29").Value = '  -0.09%  '

$ws.Range("D30").Value = '3.889.30'
$ws.Range("E30").Value = '  -0.29%  '

$ws.Range("E31").Value = '  -3.73%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.33'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.91%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.93'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.59%  '

$ws.Range("E34").Value = '  -1.85%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.29'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.14%  '

$ws.Range("D37").Value = '3.698.47'
$ws.Range("E37").Value = '  -0.46%  '

$ws.Range("E38").Value = '  -1.50%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.42'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -10.36%  '

$ws.Range("E40").Value = '  +0.73%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.997'
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.81'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.17%  '

$ws.Range("E45").Value = '  -1.72%  '

$ws.Range("B46").Value = 'Cosmos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.62'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.73%  '

$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.94'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.04%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '42.92'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +10.19%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '45.84'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.08%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '146.82'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.70%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '393.56'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.82%  '
